# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ajo" (Terminal La Palmera de La Serena)
# at rows 335-336, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 335 (pushes old rows 335.. down to 337..)
$ws.Rows("335:336").Insert()

# New row 335
$ws.Range("A335").Value = 8
$ws.Range("B335").Value = "Terminal La Palmera de La Serena"
$ws.Range("C335").Value = "Coquimbo"
$ws.Range("D335").Value = 44988
$ws.Range("E335").Value = 4
$ws.Range("F335").Value = 100112003
$ws.Range("G335").Value = "Ajo"
$ws.Range("H335").Value = "Chino"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 450
$ws.Range("K335").Value = 17000
$ws.Range("L335").Value = 18000
$ws.Range("M335").Value = 17500
$ws.Range("N335").Value = "$/caja 10 kilos"
$ws.Range("O335").Value = "China"
$ws.Range("P335").Value = 1750
$ws.Range("Q335").Value = 10
$ws.Range("R335").Value = "Hortaliza"

# New row 336
$ws.Range("A336").Value = 8
$ws.Range("B336").Value = "Terminal La Palmera de La Serena"
$ws.Range("C336").Value = "Coquimbo"
$ws.Range("D336").Value = 44988
$ws.Range("E336").Value = 4
$ws.Range("F336").Value = 100112003
$ws.Range("G336").Value = "Ajo"
$ws.Range("H336").Value = "Chino"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 400
$ws.Range("K336").Value = 20000
$ws.Range("L336").Value = 21000
$ws.Range("M336").Value = 20500
$ws.Range("N336").Value = "$/malla 10 kilos"
$ws.Range("O336").Value = "China"
$ws.Range("P336").Value = 2050
$ws.Range("Q336").Value = 10
$ws.Range("R336").Value = "Hortaliza"
